# Daily update at 8 AM UTC
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The previous "last" row (row 90) loses its special bottom-of-table
# number format (YYYY-MM-DD) and becomes a regular dated row
# (YYYY-MM-DD HH:MM:SS), matching all the other interior rows.
$ws.Range("A90").NumberFormat = "YYYY-MM-DD HH:MM:SS"

# Append the new day's data as row 91.
$ws.Range("A91").Value = 45830
$ws.Range("A91").NumberFormat = "YYYY-MM-DD"

$ws.Range("B91").Value = 384
$ws.Range("C91").Value = 390
$ws.Range("D91").Value = 391
